$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lusid_holdings")

# New "strategy" column header
$ws.Range("E1").Value = "strategy"

# Strategy values for existing rows 2-15
$strategy = @{
    2  = "Quantitative"
    3  = "Rebalance"
    4  = "Quantitative"
    5  = "Quantitative"
    6  = "Rebalance"
    7  = "Quantitative"
    8  = "Growth"
    9  = "Quantitative"
    10 = "Growth"
    11 = "Rebalance"
    12 = "Quantitative"
    13 = "Quantitative"
    14 = "Quantitative"
    15 = "Quantitative"
}

foreach ($row in $strategy.Keys) {
    $ws.Cells.Item($row, 5).Value = $strategy[$row]
}

# New holding rows 16-18
$newRows = @(
    @{ Row = 16; Isin = "JE00B4T3BW64"; Units = 1450; Currency = "GBP"; Strategy = "Rebalance" }
    @{ Row = 17; Isin = "GB0031743007"; Units = 790;  Currency = "GBP"; Strategy = "Rebalance" }
    @{ Row = 18; Isin = "GB0005603997"; Units = 2300; Currency = "GBP"; Strategy = "Growth" }
)

$holdingDate = $ws.Cells.Item(2, 1).Value2
$ws.Cells.Item(2, 1).Copy()

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $holdingDate
    $ws.Cells.Item($r, 2).Value = $entry.Isin
    $ws.Cells.Item($r, 3).Value = $entry.Units
    $ws.Cells.Item($r, 4).Value = $entry.Currency
    $ws.Cells.Item($r, 5).Value = $entry.Strategy
}

[void]$ws.Range("E18").Select()
